$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$timestamp = "2025-11-14 06:35:29"

for ($row = 2; $row -le 9; $row++) {
    $ws.Cells.Item($row, 1).Value = $timestamp
}
